$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("pythonCode")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Update cell text values ---

# A2 on the pythonCode sheet: the "sum of" python snippet
$text2 = "num1=7`nnum2=9`nsum=num1+num2`nprint(""The sum of {0} and {1} is {2}"" .format(num1,num2,sum))"
$ws1.Range("A2").Value = $text2

# A3 on the pythonCode sheet: the "sum is" python snippet
$text3 = "num1=7`nnum2=9`nsum=num1+num2`nprint(""The sum is"" sum)"
$ws1.Range("A3").Value = $text3

# Sheet2 A1 stays "pythonCode" (unchanged text, just a shared-string reindex)
$ws2.Range("A1").Value = "pythonCode"

# --- Apply wrap-text formatting to the edited cells ---
$ws1.Range("A2:A3").WrapText = $true

# --- Row heights to fit the wrapped, multi-line text ---
$ws1.Rows.Item(2).RowHeight = 135
$ws1.Rows.Item(3).RowHeight = 90

# --- Move the active selection on the pythonCode sheet to A3 ---
$ws1.Activate()
$ws1.Range("A3").Select()
